$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")
$tbl = $ws.ListObjects.Item("Table1")

# ---------------------------------------------------------------------------
# 1. Insert a new data row into Table1 at worksheet row 93 (pushes the
#    existing rows 93..133 down to 94..134, carrying the table's special
#    "last row" style down with them automatically).
# ---------------------------------------------------------------------------
$ws.Range("A93:K93").Insert(-4121)

# The newly inserted row (93) comes in with blank/default formatting; copy
# the normal data-row formatting down from the row below (now row 94, which
# used to be row 93) so it matches every other interior table row.
$ws.Range("A94:K94").Copy()
$ws.Range("A93:K93").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false

# Re-extend the table definition to include the new last row (134).
$tbl.Resize($ws.Range("A8:K134"))

# The Insert/shift mangles the calculated-column formulas on the affected
# boundary rows (new row 93, and the shifted-down former last row, now 134)
# -- restore them explicitly so they match every other row in the column.
$gFormula = '=IF(ISBLANK(Table1[[#This Row],[EARNED]]),"",Table1[[#This Row],[EARNED]])'
$ws.Range("G93").Formula = $gFormula
$ws.Range("G134").Formula = $gFormula

# ---------------------------------------------------------------------------
# 2. Fill in the new leave entry split across rows 92 and 93.
# ---------------------------------------------------------------------------
$ws.Range("B92").Value = "SL(1-0-00)"
$ws.Range("C92").Value = 1.25
$ws.Range("H92").Value = 1

$ws.Range("B93").Value = "SL(1-0-00)"
$ws.Range("H93").Value = 1

# K92 / K93 hold dates; copy the date number-format already used elsewhere
# in column K (row 86) so the new cells land on the same style as the rest
# of the date entries, then set the actual serial date values.
$ws.Range("K86").Copy()
$ws.Range("K92").PasteSpecial(-4122)
$ws.Range("K93").PasteSpecial(-4122)
$ws.Application.CutCopyMode = $false
$ws.Range("K92").Value = 45005
$ws.Range("K93").Value = 45009

# ---------------------------------------------------------------------------
# 3. Refresh the view state to mirror where the user ended up scrolled/
#    selected after the edit.
# ---------------------------------------------------------------------------
$ws.Range("B94").Select()
$excel.ActiveWindow.ScrollRow = 82
